# ---------------------------------------------------------------------------
# WeeklySummary.docx edit
#
# 1. The hidden "_GoBack" bookmark that currently sits at the top of the
#    "Wednesday July 10, 2013" paragraph is removed from there.
# 2. The "Thursday July 11, 2013 - " paragraph (the last paragraph in the
#    document) is extended with an end-dash and a new block of journal text
#    (with the same proofing marks Word's grammar/spell checker leaves
#    behind), followed by a blank paragraph and a new
#    "Friday July 12, 2013 - " paragraph.
# 3. The "_GoBack" bookmark is re-created at the end of that new
#    "Friday July 12, 2013 - " paragraph (this is where Word leaves it after
#    the last edit made to the document).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Remove the stale "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Rewrite the last paragraph ("Thursday July 11, 2013 - ") ----------
# Common run properties used throughout the journal entries.
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$pPr = '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'

# Paragraph: "Thursday July 11, 2013 ..." with the new text appended.
$thursdayPara = "<w:p>$pPr" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Thursday July 11, 2013 </w:t></w:r>" +
    "<w:r>$rPr<w:t>–</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Fixed a repaint issue when the picture was loaded. I did some more research on the JAI API for the objects that </w:t></w:r>" +
    '<w:proofErr w:type="gramStart"/>' +
    "<w:r>$rPr<w:t>are read</w:t></w:r>" +
    '<w:proofErr w:type="gramEnd"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> only like </w:t></w:r>" +
    '<w:proofErr w:type="spellStart"/>' +
    "<w:r>$rPr<w:t>PlanarImage</w:t></w:r>" +
    '<w:proofErr w:type="spellEnd"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> or Raster. </w:t></w:r>" +
    '<w:proofErr w:type="gramStart"/>' +
    "<w:r>$rPr<w:t>Started playing around with the Raster object and the pixel data.</w:t></w:r>" +
    '<w:proofErr w:type="gramEnd"/>' +
    "</w:p>"

# Blank separator paragraph.
$blankPara = "<w:p>$pPr</w:p>"

# Paragraph: "Friday July 12, 2013 - " with the relocated "_GoBack" bookmark.
$fridayPara = "<w:p>$pPr" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Friday July 12, 2013 - </w:t></w:r>" +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    "</w:p>"

$bodyXml = $thursdayPara + $blankPara + $fridayPara

$wordOpenXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    "<w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# Locate the last paragraph ("Thursday July 11, 2013 - ") and replace it.
$lastIndex = $d.Paragraphs.Count
$targetRange = $d.Paragraphs.Item($lastIndex).Range
$targetRange.InsertXML($wordOpenXml)

# InsertXML leaves a stray empty trailing paragraph behind (it cannot replace
# the implicit final paragraph mark of the story) - remove it again so the
# "Friday July 12, 2013 - " paragraph (with the bookmark) is once more the
# last paragraph in the document.
$newLastIndex = $d.Paragraphs.Count
$newFridayPara = $d.Paragraphs.Item($newLastIndex - 1)
$strayPara = $d.Paragraphs.Item($newLastIndex)
$cleanupRange = $d.Range($newFridayPara.Range.End - 1, $strayPara.Range.End)
$cleanupRange.Delete()
